# "Empezando integrador segunda entrega"
# Fill in the Precedencia (E) / Duracion (G) columns for the "Tercera entrega"
# block of activities (rows 190-268), plus a couple of incidental formatting
# tweaks the author made along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be created in this exact order so they land at
#     the same sharedStrings.xml indices the author ended up with (266-269).
$ws.Range("E198").Value = "PROY03093, PROY03094"
$ws.Range("G198").Value = 2

$ws.Range("E216").Value = "PROG03102, PROG03103"
$ws.Range("G216").Value = 5

$ws.Range("E254").Value = "FEO3122"
$ws.Range("G254").Value = 2

$ws.Range("E246").Value = "SO03114, SO03115 "
$ws.Range("G246").Value = 1

# --- Remaining Precedencia / Duracion pairs (re-use existing shared strings).
$ws.Range("E190").Value = "/"
$ws.Range("G190").Value = 1

$ws.Range("E192").Value = "/"
$ws.Range("G192").Value = 1

$ws.Range("E194").Value = "/"
$ws.Range("G194").Value = 1

$ws.Range("E196").Value = "/"
$ws.Range("G196").Value = 1

$ws.Range("E200").Value = "/"
$ws.Range("G200").Value = 3

$ws.Range("E202").Value = "/"
$ws.Range("G202").Value = 3

$ws.Range("E204").Value = "/"
$ws.Range("G204").Value = 2

$ws.Range("E206").Value = "/"
$ws.Range("G206").Value = 1

$ws.Range("E208").Value = "/"
$ws.Range("G208").Value = 2

$ws.Range("E210").Value = "/"
$ws.Range("G210").Value = 3

$ws.Range("E212").Value = "/"
$ws.Range("G212").Value = 8

# Row 214/215 (merged C214:D215) also gets shrunk to font size 10, matching
# the sibling "PROG03102" deliverable blocks.
$ws.Range("C214").Font.Size = 10
$ws.Range("D214").Font.Size = 10
$ws.Range("C215").Font.Size = 10
$ws.Range("D215").Font.Size = 10

$ws.Range("E214").Value = "PROG03102"
$ws.Range("G214").Value = 3

$ws.Range("E218").Value = "/"
$ws.Range("G218").Value = 1

$ws.Range("E220").Value = "PROG03102"
$ws.Range("G220").Value = 2

$ws.Range("E222").Value = "PROG03102"
$ws.Range("G222").Value = 4

$ws.Range("E224").Value = "PROG03102"
$ws.Range("G224").Value = 4

$ws.Range("E226").Value = "/"
$ws.Range("G226").Value = 2

$ws.Range("E228").Value = "/"
$ws.Range("G228").Value = 3

$ws.Range("E230").Value = "PROG03102"
$ws.Range("G230").Value = 3

# Row 232's Precedencia/Duracion cells pick up a top border (matching the
# plain style used elsewhere) instead of the stray highlighted one.
$ws.Range("E232").Borders.Item(8).LineStyle = 1
$ws.Range("E232").Borders.Item(8).Weight = 2
$ws.Range("F232").Borders.Item(8).LineStyle = 1
$ws.Range("F232").Borders.Item(8).Weight = 2

$ws.Range("E232").Value = "PROG03102"
$ws.Range("G232").Value = 3

$ws.Range("E234").Value = "/"
$ws.Range("G234").Value = 5

$ws.Range("E236").Value = "/"
$ws.Range("G236").Value = 3

$ws.Range("E238").Value = "SO03114"
$ws.Range("G238").Value = 3

$ws.Range("E240").Value = "SO03114"
$ws.Range("G240").Value = 2

$ws.Range("E242").Value = "SO03114"
$ws.Range("G242").Value = 2

$ws.Range("E244").Value = "SO03117"
$ws.Range("G244").Value = 3

$ws.Range("E248").Value = "/"
$ws.Range("G248").Value = 3

$ws.Range("E250").Value = "/"
$ws.Range("G250").Value = 2

$ws.Range("E252").Value = "FE03120"
$ws.Range("G252").Value = 2

$ws.Range("E256").Value = "FE03120"
$ws.Range("G256").Value = 2

$ws.Range("E258").Value = "FE03123"
$ws.Range("G258").Value = 2

$ws.Range("E260").Value = "/"
$ws.Range("G260").Value = 1

$ws.Range("E262").Value = "/"
$ws.Range("G262").Value = 2

$ws.Range("E264").Value = "/"
$ws.Range("G264").Value = 4

$ws.Range("E266").Value = "/"
$ws.Range("G266").Value = 3

$ws.Range("E268").Value = "/"
$ws.Range("G268").Value = 2

# --- Scroll position / selection, matching where the author was working.
$ws.Range("G208:H209").Select()
